$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.226.32'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.603.69'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.51'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.44'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("E10").Value = '  +1.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.336'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.060.43'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.161.26'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.58'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.622.93'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.42'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.13'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.21'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0756'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.06%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.73'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.67'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.99'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '37.13'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.36%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.834'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.57'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '274.93'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0524'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.954.63'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0224'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.53'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.02'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.84%  '
